$wb = $excel.ActiveWorkbook
$ingreso = $wb.Worksheets.Item("Ingreso")
$gastos = $wb.Worksheets.Item("Gastos")

# --- Fix up the tail of the existing "Ingreso" data (rows 586-588) ---
$ingreso.Range("D586").Value = "Aporte"
$ingreso.Range("D587").Value = "Aporte"
$ingreso.Range("B588").Value = "kukito"
$ingreso.Range("D588").Value = "Aporte"

# --- Append new contribution rows (week of 2024-01-14, serial 45305) ---
$week1 = @(
    @("Kukito", 0),
    @("Kibelo", 0),
    @("Rubio", 0),
    @("Jordan", 0),
    @("Johan", 0),
    @("Chamo", 100),
    @("Joel", 0),
    @("Julio", 0),
    @("Javier", 0),
    @("Yeyo", 0)
)

$r = 589
foreach ($entry in $week1) {
    $ingreso.Cells.Item($r, 1).Value = 45305
    $ingreso.Cells.Item($r, 2).Value = $entry[0]
    $ingreso.Cells.Item($r, 3).Value = $entry[1]
    $ingreso.Cells.Item($r, 4).Value = "Aporte"
    $r = $r + 1
}

# --- One late entry still dated 2024-01-06 (serial 45297) ---
$ingreso.Cells.Item(599, 1).Value = 45297
$ingreso.Cells.Item(599, 2).Value = "Joel"
$ingreso.Cells.Item(599, 3).Value = 0
$ingreso.Cells.Item(599, 4).Value = "Aporte"

# --- Week of 2024-01-21 (serial 45312) ---
# Some of these rows were pasted in without the usual "amount" number
# style, so column C keeps the worksheet's bare default format there.
$week2 = @(
    @("Michy", 100, $false),
    @("Frandy", 100, $false),
    @("Jordan", 0, $true),
    @("Johan", 0, $true),
    @("Chamo", 0, $false),
    @("Orlando", 100, $false),
    @("Kukito", 0, $false),
    @("Rubio", 0, $false),
    @("Kibelo", 0, $false),
    @("Joel", 1000, $false),
    @("Rayder", 100, $false),
    @("Yeyo", 0, $false)
)

$r = 600
foreach ($entry in $week2) {
    $ingreso.Cells.Item($r, 1).Value = 45312
    $ingreso.Cells.Item($r, 2).Value = $entry[0]
    $ingreso.Cells.Item($r, 3).Value = $entry[1]
    if (-not $entry[2]) {
        $ingreso.Cells.Item($r, 3).Style = "Normal"
    }
    $ingreso.Cells.Item($r, 4).Value = "Aporte"
    $r = $r + 1
}

# --- Bump the conditional-format rule's priority (Excel re-numbers it as the
#     duplicate-values highlight gets touched while the new rows are added) ---
$dupRule = $ingreso.Range("B1:B1048576").FormatConditions.Item(1)
$dupRule.Priority = 4

# --- "Gastos" sheet: new expense entry for 2024-01-21 ---
$gastos.Cells.Item(77, 1).Value = 45312
$gastos.Cells.Item(77, 2).Value = "Agua"
$gastos.Cells.Item(77, 3).Value = 140

$gastos.Activate()
$gastos.Range("A77").Select()

# --- Leave "Ingreso" as the active tab with the cursor where data entry
#     stopped, matching the author's final on-screen state ---
$ingreso.Activate()
$ingreso.Range("C589").Select()
